$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price and Volume columns remain text (avoid Excel auto-numeric conversion)
$ws.Columns.Item(4).NumberFormat = "@"
$ws.Columns.Item(5).NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '37.830.27'
$ws.Range("E2").Value = '  +0.03%  '

# Row 3
$ws.Range("D3").Value = '2.076.51'
$ws.Range("E3").Value = '  -1.02%  '

# Row 4
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").Value = '233.11'
$ws.Range("E5").Value = '  +0.41%  '

# Row 6
$ws.Range("D6").Value = '0.625'
$ws.Range("E6").Value = '  +0.27%  '

# Row 7
$ws.Range("D7").Value = '59.10'
$ws.Range("E7").Value = '  +1.89%  '

# Row 8
$ws.Range("E8").Value = '  -0.04%  '

# Row 9
$ws.Range("D9").Value = '0.394'
$ws.Range("E9").Value = '  +1.25%  '

# Row 10
$ws.Range("D10").Value = '0.0789'
$ws.Range("E10").Value = '  +1.07%  '

# Row 11
$ws.Range("E11").Value = '  +1.91%  '

# Row 12
$ws.Range("B12").Value = 'Chainlink'
$ws.Range("C12").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D12").Value = '14.76'
$ws.Range("E12").Value = '  +1.79%  '

# Row 13
$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").Value = '21.20'
$ws.Range("E13").Value = '  +0.27%  '

# Row 14
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = '0.772'
$ws.Range("E14").Value = '  +0.54%  '

# Row 15
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '5.35'
$ws.Range("E15").Value = '  +2.30%  '

# Row 16
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '2.056.11'
$ws.Range("E16").Value = '  -1.97%  '

# Row 17
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '37.705.95'
$ws.Range("E17").Value = '  -0.21%  '

# Row 18
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").Value = '6.14'
$ws.Range("E18").Value = '  -0.34%  '

# Row 19
$ws.Range("B19").Value = 'Litecoin'
$ws.Range("C19").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D19").Value = '71.54'
$ws.Range("E19").Value = '  +1.24%  '

# Row 20
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.0₃0850'
$ws.Range("E20").Value = '  +3.33%  '

# Row 21
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = '228.22'
$ws.Range("E21").Value = '  +0.13%  '

# Row 22
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.04%  '

# Row 23
$ws.Range("B23").Value = 'PancakeSwap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D23").Value = '2.41'
$ws.Range("E23").Value = '  +1.17%  '

# Row 24
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").Value = '2.37'
$ws.Range("E24").Value = '  -1.23%  '

# Row 25
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = '170.58'
$ws.Range("E25").Value = '  +1.48%  '

# Row 26
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '9.19'
$ws.Range("E26").Value = '  +2.90%  '

# Row 27
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").Value = '0.134'
$ws.Range("E27").Value = '  -4.45%  '

# Row 28
$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").Value = '1.42'
$ws.Range("E28").Value = '  -0.60%  '

# Row 29
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '19.48'
$ws.Range("E29").Value = '  +0.12%  '

# Row 30
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = '0.121'
$ws.Range("E30").Value = '  +1.61%  '

# Row 31
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '4.73'
$ws.Range("E31").Value = '  +2.22%  '

# Row 32
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = '4.74'
$ws.Range("E32").Value = '  +3.47%  '

# Row 33
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '0.0632'
$ws.Range("E33").Value = '  +0.95%  '

# Row 34
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").Value = '2.49'
$ws.Range("E34").Value = '  -0.54%  '

# Row 35
$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").Value = '1.82'
$ws.Range("E35").Value = '  -0.31%  '

# Row 36
$ws.Range("D36").Value = '3.40'
$ws.Range("E36").Value = '  -0.48%  '

# Row 37
$ws.Range("B37").Value = 'BinanceUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.00%  '

# Row 38
$ws.Range("B38").Value = 'THORChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D38").Value = '5.41'
$ws.Range("E38").Value = '  +0.28%  '

# Row 39
$ws.Range("B39").Value = 'Cronos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D39").Value = '0.0984'
$ws.Range("E39").Value = '  -1.37%  '

# Row 40
$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D40").Value = '99.22'
$ws.Range("E40").Value = '  +1.34%  '

# Row 41
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.0216'
$ws.Range("E41").Value = '  +0.82%  '

# Row 42
$ws.Range("B42").Value = 'HuobiToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D42").Value = '2.88'
$ws.Range("E42").Value = '  -1.81%  '

# Row 43
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").Value = '16.70'
$ws.Range("E43").Value = '  +6.81%  '

# Row 44
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '1.443.08'
$ws.Range("E44").Value = '  -0.90%  '

# Row 45
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").Value = '1.16'
$ws.Range("E45").Value = '  -0.77%  '

# Row 46
$ws.Range("E46").Value = '  +5.81%  '

# Row 47
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").Value = '1.06'
$ws.Range("E47").Value = '  +0.37%  '

# Row 48
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").Value = '7.39'
$ws.Range("E48").Value = '  +0.52%  '

# Row 49
$ws.Range("B49").Value = 'MXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D49").Value = '3.00'
$ws.Range("E49").Value = '  -0.47%  '

# Row 50
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.267.01'
$ws.Range("E50").Value = '  -1.13%  '

# Row 51
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").Value = '46.64'
$ws.Range("E51").Value = '  +0.77%  '
